$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.989.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.847.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.884.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.430.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.70%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +8.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "313.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.61%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0962"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0520"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.411"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.57%  "
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  -0.32%  "
